$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Delete row 548 ("「早起きする方法」" entry) entirely; subsequent rows shift up.
$ws.Rows.Item(548).Delete()
